# Update mods data [2026-01-23 15:14:22]
# Append a new daily data row (2026/01/23, 逃离鸭科夫, 1156) to the bottom
# of the ModCounts sheet, matching the look/style of the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last populated row (currently row 73, A1:C73) and append after it.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1
$newRow = $lastRow + 1

# Copy the formatting (centered alignment style) of the last data row down
# onto the new row, so the new cells pick up the same style as the rest of
# the table (instead of the workbook's default/general style).
$ws.Range("A" + $lastRow + ":C" + $lastRow).Copy()
$ws.Range("A" + $newRow + ":C" + $newRow).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Column A holds a date written as plain text (e.g. "2026/01/22"), not a
# real Excel date serial. Mark the new cell as Text first so typing a
# slash-separated, date-shaped string doesn't get auto-converted into a
# date value.
$ws.Range("A" + $newRow).NumberFormat = "@"
$ws.Range("A" + $newRow).Value = "2026/01/23"

$ws.Range("B" + $newRow).Value = "逃离鸭科夫"
$ws.Range("C" + $newRow).Value = 1156
